# Update "想去人数" (F column) counts on sheets "展览" and "全部类型"
# to reflect the newly generated numbers (commit: output generated at 456a3b4)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    2  = 1044
    8  = 1654
    9  = 6076
    12 = 277
    16 = 5409
    18 = 1264
    23 = 256
    28 = 378
    29 = 71
    35 = 59
    36 = 61
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
